$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values ---
# Write order chosen to reproduce the shared-string table ordering of the
# target workbook: new strings are introduced in this exact sequence.
$ws.Range("E11").Value = "BC,VEG"
$ws.Range("B2").Value  = "Eggs / Bacon / Tater Tots /  Salsa / Nacho Cheese"
$ws.Range("C2").Value  = "Wheat, milk, sulphites, egg."
$ws.Range("D2").Value  = "Castle Cheese Nacho Blend"
$ws.Range("B3").Value  = "Eggs / Pork Sausage / Tater Tots /  Salsa / Nacho Cheese"
$ws.Range("D3").Value  = " Castle Cheese Nacho Blend"
$ws.Range("B4").Value  = "Eggs / Spinach / Tater Tots /  Salsa / Nacho Cheese"
$ws.Range("E8").Value  = "BC,VGN,DF"
$ws.Range("D8").Value  = "Maureen's Tofu Roll"
$ws.Range("B9").Value  = "Flaky Pastry / Ground Beef Sausage / Spices"
$ws.Range("B10").Value = "Flaky Pastry / Ground Seasoned Chicken Sausage / Spices"
$ws.Range("B8").Value  = "Flaky Pastry / Tofu Scramble / Spices"
$ws.Range("E9").Value  = "DF"
$ws.Range("F8").Value  = "Tofu_Scramble_Roll"
$ws.Range("F9").Value  = "Beef_Sausage_Roll"
$ws.Range("F2").Value  = "Bacon_Breakfast_Burrito"
$ws.Range("F4").Value  = "Falafel_Wrap"

# Remaining changed cells that only reuse strings introduced above
$ws.Range("C3").Value  = "Wheat, milk, sulphites, egg."
$ws.Range("F3").Value  = "Bacon_Breakfast_Burrito"
$ws.Range("C4").Value  = "Wheat, milk, sulphites, egg."
$ws.Range("D4").Value  = " Castle Cheese Nacho Blend"
$ws.Range("C9").Value  = "Wheat, soy."
$ws.Range("C10").Value = "Wheat, soy."
$ws.Range("E10").Value = "DF"
$ws.Range("F10").Value = "Beef_Sausage_Roll"
$ws.Range("F11").Value = "Egg_and_Cheese_Breakfast_Sandwich"

# --- Style changes: clear the highlighted/wrapped formatting on cells that  ---
# --- no longer use it (reset to the workbook's default "Normal" style)     ---
$ws.Range("F2").Style  = "Normal"
$ws.Range("B3").Style  = "Normal"
$ws.Range("F3").Style  = "Normal"
$ws.Range("F4").Style  = "Normal"
$ws.Range("F11").Style = "Normal"

# --- Row height: row 4 loses its custom height and reverts to the default ---
$ws.Rows.Item(4).AutoFit()

# --- Selection moves from F18 to F10 ---
$ws.Range("F10").Select()
